# Daily attendance processing - 2025-12-02 19:05:43
# Normalise the "Recorded By" column (G): when a cell lists multiple
# recorders separated by ", ", move the first-listed name to the end
# of the list (rotate left by one) instead of leading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Text

    if ($val -eq $null -or $val -eq "") {
        continue
    }

    $parts = $val -split ", "

    if ($parts.Length -gt 1) {
        $rotated = $parts[1..($parts.Length - 1)] + $parts[0]
        $newVal = $rotated -join ", "
        $cell.Value = $newVal
    }
}
